# Update "想去人数" (column F) values on the "展览" and "全部类型" sheets
# to match newly scraped counts.

$wb = $excel.ActiveWorkbook

$updates = @{
    2  = 70
    3  = 8580
    4  = 1533
    7  = 275
    10 = 131
    11 = 46
    13 = 1283
    14 = 302
    18 = 141
    19 = 84
    20 = 128
    21 = 116
}

$sheetNames = @("展览", "全部类型")

foreach ($sheetName in $sheetNames) {
    $ws = $wb.Worksheets.Item($sheetName)
    foreach ($row in $updates.Keys) {
        $ws.Range("F$row").Value = $updates[$row]
    }
}

$wb.Save()
